$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column in H1, matching the header styling
# already used by the other header cells (bold/bordered via column G's style).
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Seed the new "Save" column with 0 for every existing data row (rows 2-9).
$ws.Range("H2:H9").Value = 0
